$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new order rows (15-18, written into sheet rows 16-19) ---

# Row 16: Easter / Toy / Mr. Hoppers_2 / E4835T
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Easter"
$ws.Range("C16").Value = "Toy"
$ws.Range("D16").Value = "Mr. Hoppers_2"
$ws.Range("E16").Value = 18
$ws.Range("F16").Value = "E4835T"
$ws.Range("G16").Value = "Learn the alphabet with the interactive toys for infants."
$ws.Range("H16").Value = "Y"
$ws.Range("I16").Value = 1
$ws.Range("P16").Value = 30
$ws.Range("Q16").Value = "Pink"

# Row 17: duplicate of row 16
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Easter"
$ws.Range("C17").Value = "Toy"
$ws.Range("D17").Value = "Mr. Hoppers_2"
$ws.Range("E17").Value = 18
$ws.Range("F17").Value = "E4835T"
$ws.Range("G17").Value = "Learn the alphabet with the interactive toys for infants."
$ws.Range("H17").Value = "Y"
$ws.Range("I17").Value = 1
$ws.Range("P17").Value = 30
$ws.Range("Q17").Value = "Pink"

# Row 18: Christmas / Toy / Santas Workshop - Essentials Edition / C1230T
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Christmas"
$ws.Range("C18").Value = "Toy"
$ws.Range("D18").Value = "Santas Workshop - Essentials Edition"
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = "C1230T"
$ws.Range("G18").Value = "The most sought after christmas present! Get yours today!"
$ws.Range("H18").Value = "N"
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = "50,90"
$ws.Range("K18").Value = 4

# Row 19: duplicate of row 18, plus an extra jump_height value
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Christmas"
$ws.Range("C19").Value = "Toy"
$ws.Range("D19").Value = "Santas Workshop - Essentials Edition"
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = "C1230T"
$ws.Range("G19").Value = "The most sought after christmas present! Get yours today!"
$ws.Range("H19").Value = "N"
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = "50,90"
$ws.Range("K19").Value = 4
$ws.Range("M19").Value = 10

# --- View / window adjustments ---
$ws.Range("G16:G19").WrapText = $true
$ws.Rows.Item(16).RowHeight = 16
$ws.Rows.Item(17).RowHeight = 16
$ws.Rows.Item(18).RowHeight = 16
$ws.Rows.Item(19).RowHeight = 16

$excel.ActiveWindow.Zoom = 165
$ws.Range("K19").Select()

$excel.ActiveWindow.WindowState = -4137
$excel.Left = -25600
$excel.Top = 0
$excel.Width = 25600
$excel.Height = 28800
